$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("PDiCCpDoC")

# Update onshore wind, solar PV and offshore wind percentages, highlighting
# the changed cells in yellow.
$wsData.Range("B6").Value = 0.15
$wsData.Range("B6").Interior.Color = 65535

$wsData.Range("B7").Value = 0.6
$wsData.Range("B7").Interior.Color = 65535

$wsData.Range("B14").Value = 0.15
$wsData.Range("B14").Interior.Color = 65535

# Update the active sheet/selection so the PDiCCpDoC sheet is the one shown
# when the workbook is reopened.
$wsAbout.Range("B17").Select()
$wsData.Activate()
$wsData.Range("B7").Select()
